# Remove the two ".NET STANDARD" API-count slides (old slide 3 and slide 4)
# from the deck. These were the slides with id="452" (13k/32k APIs bars) and
# id="438" (13k APIs with rectangle/arrow). Their notes slides
# (notesSlide1.xml / notesSlide2.xml) are removed automatically along with
# them since nothing else references them.
#
# Deleting index 3 twice removes both slides: after the first delete, the
# slide that used to be 4th shifts into position 3.
$p = $ppt.ActivePresentation
$p.Slides.Item(3).Delete()
$p.Slides.Item(3).Delete()
